$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new "orig-code detecting result" columns (F, G) for data rows 14-16 ---
$ws.Range("F14").Value = 0.733
$ws.Range("G14").Value = 0.576

$ws.Range("F15").Value = 0.739
$ws.Range("G15").Value = 0.62

$ws.Range("F16").Value = 0.748
$ws.Range("G16").Value = 0.598

# --- Add the new summary row 18: average of the last three rows (14:16) for each column ---
$ws.Range("B18").Formula = "=SUM(B14:B16)/3"
$ws.Range("C18:F18").Formula = "=SUM(C14:C16)/3"
$ws.Range("G18").Formula = "=SUM(G14:G16)/3"

# --- Reposition/resize the chart to make room for the new data (user dragged it) ---
$co = $ws.ChartObjects().Item(1)
$co.Left = 460.7054133858268
$co.Top = 93.89992125984251
$co.Width = 385.0625
$co.Height = 216.1071653543307

# --- Final selection left on the sheet after the edits ---
[void]$ws.Range("O19").Select()
